$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Imad"
$ws.Range("C2").Value = "Shehadeh"
$ws.Range("D2").Formula = "=TEXT(25874125,""0"")"
$ws.Range("D2").Copy()
$ws.Range("D2").PasteSpecial(-4163)
$ws.Range("E2").Value = "VIP"
